# "not yet contacted list"
# Remove the "Foulen El Fouleni" record (row 2) and add a new "Mohamed akram
# bouzaiane" record, shifting the remaining rows so the sheet reads:
#   2: bouden eya
#   3: Feriel chouaieb
#   4: Med Amine Allani
#   5: Mohamed akram bouzaiane   (new)
#   6: Rani ZOUAOUI
#
# Column A keeps its original literal 0..4 index values, so we rewrite every
# data cell explicitly instead of relying on EntireRow.Delete/Insert (which
# would also shift the literal numbers stored in column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Drop every existing hyperlink; we'll re-add the ones we still need. ----
$ws.Range("E2").Hyperlinks.Delete()
$ws.Range("E3").Hyperlinks.Delete()
$ws.Range("E4").Hyperlinks.Delete()
$ws.Range("E5").Hyperlinks.Delete()
$ws.Range("E6").Hyperlinks.Delete()

# ---- Wipe the whole data block (keeps per-cell styles as empty stubs). ----
$ws.Range("A2:L6").ClearContents()

# E5 had the leftover "hyperlink" style stub from the old row; the new row 5
# (Mohamed) has no facebook-link cell at all, so drop it completely (style
# and all) rather than leaving an empty styled cell behind.
$ws.Range("E5").Clear()

function Set-TextValue($addr, $text) {
    # Force text storage so numeric-looking strings (phone numbers, etc.)
    # keep their trailing ".0" instead of being coerced into a number.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# ---------------------------------------------------------------------
# Row 2 : bouden eya
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 0
Set-TextValue "B2" "bouden eya"
Set-TextValue "C2" "94252435.0"
Set-TextValue "D2" "eyabou22@gmail.com"
Set-TextValue "E2" "https://www.facebook.com/eya22T/"
Set-TextValue "F2" "done"
Set-TextValue "H2" "ilyes"
Set-TextValue "I2" "no"

# ---------------------------------------------------------------------
# Row 3 : Feriel chouaieb
# ---------------------------------------------------------------------
$ws.Range("A3").Value = 1
Set-TextValue "B3" "Feriel chouaieb"
Set-TextValue "C3" "54234540.0"
Set-TextValue "D3" "Ferielchouaieb@gmail.com"
Set-TextValue "E3" "feriel chouaieb"
Set-TextValue "F3" "done"
Set-TextValue "G3" "local"
Set-TextValue "I3" "no"
Set-TextValue "K3" "2002-05-27 00:00:00"

# ---------------------------------------------------------------------
# Row 4 : Med Amine Allani
# ---------------------------------------------------------------------
$ws.Range("A4").Value = 2
Set-TextValue "B4" "Med Amine Allani"
Set-TextValue "C4" "53265593.0"
Set-TextValue "D4" "medamineallani@gmail.com"
Set-TextValue "E4" "https://www.facebook.com/minouallani"
Set-TextValue "F4" "done"
Set-TextValue "G4" "local"
Set-TextValue "H4" "mariem"
Set-TextValue "I4" "BIG NO"

# ---------------------------------------------------------------------
# Row 5 : Mohamed akram bouzaiane (new applicant)
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 3
Set-TextValue "B5" "Mohamed akram bouzaiane"
Set-TextValue "C5" "27910734.0"
Set-TextValue "D5" "mohamedakrambouzaiane@gmail.com"
Set-TextValue "F5" "done"
Set-TextValue "G5" "salle"
Set-TextValue "H5" "mariem najjar+fatma"
Set-TextValue "I5" "NO"

# ---------------------------------------------------------------------
# Row 6 : Rani ZOUAOUI
# ---------------------------------------------------------------------
$ws.Range("A6").Value = 4
Set-TextValue "B6" "Rani ZOUAOUI"
Set-TextValue "C6" "99548790.0"
Set-TextValue "D6" "ranizouaouicontact@gmail.com"
Set-TextValue "E6" "https://www.facebook.com/rani.zouaoui.775"
Set-TextValue "F6" "done"
Set-TextValue "G6" "local"
Set-TextValue "H6" "malek bokri"
Set-TextValue "I6" "NO"

# ---- Re-add hyperlinks for the facebook-link cells that still have one. ----
$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.facebook.com/eya22T/")
$ws.Range("E2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E4"), "https://www.facebook.com/minouallani")
$ws.Range("E4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E6"), "https://www.facebook.com/rani.zouaoui.775")
$ws.Range("E6").Style = "Hyperlink"
